# Change "ASB" (Azure Security Baseline) to "MCSB" (Microsoft Cloud Security
# Benchmark) in the header labels of both worksheets.
# Commit message: "Change ASB to MCSB (#117)"

$wb = $excel.ActiveWorkbook

# --- Sheet "Contoso Controls" ---
$ws1 = $wb.Worksheets.Item("Contoso Controls")
$ws1.Range("G3").Value = "MCSB Mapping"
$ws1.Range("H3").Value = "MCSB Guidance"
$ws1.Range("I3").Value = "MCSB Policy"

# --- Sheet "Contoso Controls (Complete)" ---
$ws2 = $wb.Worksheets.Item("Contoso Controls (Complete)")
$ws2.Range("F3").Value = "MCSB Mapping"
$ws2.Range("G3").Value = "MCSB Guidance"
$ws2.Range("H3").Value = "MCSB Policy"
